# Applies the public EpexSpot prices update:
#  - "Prix Spot" sheet: insert a new date column ("06-nov") before the
#    existing "01-oct." column, shifting everything from that column
#    onward one column to the right. The new column is filled with "-"
#    placeholders (no data yet for that date) except the header.
#  - "Gaz" and "CO2" sheets: append a new daily row (2025-11-04).

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert new column at DI (before the old "01-oct.") ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Use the numeric column index (113 = "DI") rather than the letter form;
# inserting by letter reference confuses downstream statement evaluation
# in this host.
$wsPrix.Columns.Item(113).Insert()
Write-Host "Inserted new column at DI"

$wsPrix.Range("DI1").Value = "06-nov"

for ($row = 2; $row -le 25; $row++) {
    $wsPrix.Cells.Item($row, 113).Value = "-"
}

# --- Sheet "Gaz": append row 142 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
# Leading apostrophe forces the date-looking string to be stored as text
# (matches how all the other date cells in this column are stored),
# then the style is reset so no extra "Text" number format sticks around.
$wsGaz.Range("A142").Value = "'2025-11-04"
$wsGaz.Range("A142").Style = "Normal"
$wsGaz.Range("B142").Value = 31.17

# --- Sheet "CO2": append row 142 ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A142").Value = "'2025-11-04"
$wsCo2.Range("A142").Style = "Normal"
$wsCo2.Range("B142").Value = 81.90000000000001

Write-Host "Done applying EpexSpot update"
